$wb = $excel.ActiveWorkbook

# --- Step 1: repurpose the existing '总计' sheet (keeps its sheetId/rId) as '2022-Q1' ---
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# --- Step 2: make a formatted copy of an existing quarter sheet to become the refreshed '总计' sheet ---
$quarterRef = $wb.Worksheets.Item("2021-Q4")
$quarterRef.Copy([Type]::Missing, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"
$total.Cells.Clear()

# --- copy header / index-column formatting onto the new 2022-Q1 sheet ---
$quarterRef.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$quarterRef.Range("A2").Copy()
$q1.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2022-Q1 header text ---
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- 2022-Q1 data rows ---
$q1.Range("A2").Value = 0
$c = $q1.Range("B2"); $c.NumberFormat = "@"; $c.Value = "011136"; $c.Style = "Normal"
$q1.Range("C2").Value = "广发盛兴混合A"
$c = $q1.Range("D2"); $c.NumberFormat = "@"; $c.Value = "22.19"; $c.Style = "Normal"
$c = $q1.Range("E2"); $c.NumberFormat = "@"; $c.Value = "93.21"; $c.Style = "Normal"
$c = $q1.Range("F2"); $c.NumberFormat = "@"; $c.Value = "5.57"; $c.Style = "Normal"
$c = $q1.Range("G2"); $c.NumberFormat = "@"; $c.Value = "1.2360"; $c.Style = "Normal"
$q1.Range("H2").Value = 7

$q1.Range("A3").Value = 1
$c = $q1.Range("B3"); $c.NumberFormat = "@"; $c.Value = "270021"; $c.Style = "Normal"
$q1.Range("C3").Value = "广发聚瑞混合A"
$c = $q1.Range("D3"); $c.NumberFormat = "@"; $c.Value = "24.67"; $c.Style = "Normal"
$c = $q1.Range("E3"); $c.NumberFormat = "@"; $c.Value = "93.75"; $c.Style = "Normal"
$c = $q1.Range("F3"); $c.NumberFormat = "@"; $c.Value = "4.98"; $c.Style = "Normal"
$c = $q1.Range("G3"); $c.NumberFormat = "@"; $c.Value = "1.2286"; $c.Style = "Normal"
$q1.Range("H3").Value = 7

$q1.Range("A4").Value = 2
$c = $q1.Range("B4"); $c.NumberFormat = "@"; $c.Value = "010161"; $c.Style = "Normal"
$q1.Range("C4").Value = "广发瑞安精选股票A"
$c = $q1.Range("D4"); $c.NumberFormat = "@"; $c.Value = "8.03"; $c.Style = "Normal"
$c = $q1.Range("E4"); $c.NumberFormat = "@"; $c.Value = "92.30"; $c.Style = "Normal"
$c = $q1.Range("F4"); $c.NumberFormat = "@"; $c.Value = "6.31"; $c.Style = "Normal"
$c = $q1.Range("G4"); $c.NumberFormat = "@"; $c.Value = "0.5067"; $c.Style = "Normal"
$q1.Range("H4").Value = 6

$q1.Range("A5").Value = 3
$c = $q1.Range("B5"); $c.NumberFormat = "@"; $c.Value = "159883"; $c.Style = "Normal"
$q1.Range("C5").Value = "永赢中证全指医疗器械交易型开放式指数证券投资基金"
$c = $q1.Range("D5"); $c.NumberFormat = "@"; $c.Value = "9.17"; $c.Style = "Normal"
$c = $q1.Range("E5"); $c.NumberFormat = "@"; $c.Value = "99.23"; $c.Style = "Normal"
$c = $q1.Range("F5"); $c.NumberFormat = "@"; $c.Value = "2.64"; $c.Style = "Normal"
$c = $q1.Range("G5"); $c.NumberFormat = "@"; $c.Value = "0.2421"; $c.Style = "Normal"
$q1.Range("H5").Value = 10

$q1.Range("A6").Value = 4
$c = $q1.Range("B6"); $c.NumberFormat = "@"; $c.Value = "001305"; $c.Style = "Normal"
$q1.Range("C6").Value = "九泰天富改革新动力混合A"
$c = $q1.Range("D6"); $c.NumberFormat = "@"; $c.Value = "3.74"; $c.Style = "Normal"
$c = $q1.Range("E6"); $c.NumberFormat = "@"; $c.Value = "88.86"; $c.Style = "Normal"
$c = $q1.Range("F6"); $c.NumberFormat = "@"; $c.Value = "5.02"; $c.Style = "Normal"
$c = $q1.Range("G6"); $c.NumberFormat = "@"; $c.Value = "0.1877"; $c.Style = "Normal"
$q1.Range("H6").Value = 10

$q1.Range("A7").Value = 5
$c = $q1.Range("B7"); $c.NumberFormat = "@"; $c.Value = "011137"; $c.Style = "Normal"
$q1.Range("C7").Value = "广发盛兴混合C"
$c = $q1.Range("D7"); $c.NumberFormat = "@"; $c.Value = "2.10"; $c.Style = "Normal"
$c = $q1.Range("E7"); $c.NumberFormat = "@"; $c.Value = "93.21"; $c.Style = "Normal"
$c = $q1.Range("F7"); $c.NumberFormat = "@"; $c.Value = "5.57"; $c.Style = "Normal"
$c = $q1.Range("G7"); $c.NumberFormat = "@"; $c.Value = "0.1170"; $c.Style = "Normal"
$q1.Range("H7").Value = 7

$q1.Range("A8").Value = 6
$c = $q1.Range("B8"); $c.NumberFormat = "@"; $c.Value = "010026"; $c.Style = "Normal"
$q1.Range("C8").Value = "广发聚瑞混合C"
$c = $q1.Range("D8"); $c.NumberFormat = "@"; $c.Value = "1.07"; $c.Style = "Normal"
$c = $q1.Range("E8"); $c.NumberFormat = "@"; $c.Value = "93.75"; $c.Style = "Normal"
$c = $q1.Range("F8"); $c.NumberFormat = "@"; $c.Value = "4.98"; $c.Style = "Normal"
$c = $q1.Range("G8"); $c.NumberFormat = "@"; $c.Value = "0.0533"; $c.Style = "Normal"
$q1.Range("H8").Value = 7

$q1.Range("A9").Value = 7
$c = $q1.Range("B9"); $c.NumberFormat = "@"; $c.Value = "010162"; $c.Style = "Normal"
$q1.Range("C9").Value = "广发瑞安精选股票C"
$c = $q1.Range("D9"); $c.NumberFormat = "@"; $c.Value = "0.65"; $c.Style = "Normal"
$c = $q1.Range("E9"); $c.NumberFormat = "@"; $c.Value = "92.30"; $c.Style = "Normal"
$c = $q1.Range("F9"); $c.NumberFormat = "@"; $c.Value = "6.31"; $c.Style = "Normal"
$c = $q1.Range("G9"); $c.NumberFormat = "@"; $c.Value = "0.0410"; $c.Style = "Normal"
$q1.Range("H9").Value = 6

$q1.Range("A10").Value = 8
$c = $q1.Range("B10"); $c.NumberFormat = "@"; $c.Value = "009912"; $c.Style = "Normal"
$q1.Range("C10").Value = "九泰天富改革新动力混合C"
$c = $q1.Range("D10"); $c.NumberFormat = "@"; $c.Value = "0.59"; $c.Style = "Normal"
$c = $q1.Range("E10"); $c.NumberFormat = "@"; $c.Value = "88.86"; $c.Style = "Normal"
$c = $q1.Range("F10"); $c.NumberFormat = "@"; $c.Value = "5.02"; $c.Style = "Normal"
$c = $q1.Range("G10"); $c.NumberFormat = "@"; $c.Value = "0.0296"; $c.Style = "Normal"
$q1.Range("H10").Value = 10

$q1.Range("A11").Value = 9
$c = $q1.Range("B11"); $c.NumberFormat = "@"; $c.Value = "008437"; $c.Style = "Normal"
$q1.Range("C11").Value = "九泰行业优选灵活配置混合A"
$c = $q1.Range("D11"); $c.NumberFormat = "@"; $c.Value = "0.11"; $c.Style = "Normal"
$c = $q1.Range("E11"); $c.NumberFormat = "@"; $c.Value = "51.13"; $c.Style = "Normal"
$c = $q1.Range("F11"); $c.NumberFormat = "@"; $c.Value = "4.08"; $c.Style = "Normal"
$c = $q1.Range("G11"); $c.NumberFormat = "@"; $c.Value = "0.0045"; $c.Style = "Normal"
$q1.Range("H11").Value = 8

$q1.Range("A12").Value = 10
$c = $q1.Range("B12"); $c.NumberFormat = "@"; $c.Value = "008438"; $c.Style = "Normal"
$q1.Range("C12").Value = "九泰行业优选灵活配置混合C"
$c = $q1.Range("D12"); $c.NumberFormat = "@"; $c.Value = "0.06"; $c.Style = "Normal"
$c = $q1.Range("E12"); $c.NumberFormat = "@"; $c.Value = "51.13"; $c.Style = "Normal"
$c = $q1.Range("F12"); $c.NumberFormat = "@"; $c.Value = "4.08"; $c.Style = "Normal"
$c = $q1.Range("G12"); $c.NumberFormat = "@"; $c.Value = "0.0024"; $c.Style = "Normal"
$q1.Range("H12").Value = 8

$q1.Range("A13").Value = 11
$c = $q1.Range("B13"); $c.NumberFormat = "@"; $c.Value = "010999"; $c.Style = "Normal"
$q1.Range("C13").Value = "兴华瑞丰混合A"
$c = $q1.Range("D13"); $c.NumberFormat = "@"; $c.Value = "0.06"; $c.Style = "Normal"
$c = $q1.Range("E13"); $c.NumberFormat = "@"; $c.Value = "29.21"; $c.Style = "Normal"
$c = $q1.Range("F13"); $c.NumberFormat = "@"; $c.Value = "3.78"; $c.Style = "Normal"
$c = $q1.Range("G13"); $c.NumberFormat = "@"; $c.Value = "0.0023"; $c.Style = "Normal"
$q1.Range("H13").Value = 3

$q1.Range("A14").Value = 12
$c = $q1.Range("B14"); $c.NumberFormat = "@"; $c.Value = "011000"; $c.Style = "Normal"
$q1.Range("C14").Value = "兴华瑞丰混合C"
$c = $q1.Range("D14"); $c.NumberFormat = "@"; $c.Value = "0.05"; $c.Style = "Normal"
$c = $q1.Range("E14"); $c.NumberFormat = "@"; $c.Value = "29.21"; $c.Style = "Normal"
$c = $q1.Range("F14"); $c.NumberFormat = "@"; $c.Value = "3.78"; $c.Style = "Normal"
$c = $q1.Range("G14"); $c.NumberFormat = "@"; $c.Value = "0.0019"; $c.Style = "Normal"
$q1.Range("H14").Value = 3

# --- copy header / index-column formatting onto the refreshed 总计 sheet ---
$quarterRef.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$quarterRef.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 总计 header text ---
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# --- 总计 data rows (2022-Q1 inserted at the top, remaining quarters shifted down) ---
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 13
$total.Range("D2").Value = 3.65

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 11
$total.Range("D3").Value = 6.71

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 12
$total.Range("D4").Value = 7.47

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 20
$total.Range("D5").Value = 11.45

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 52
$total.Range("D6").Value = 37.41

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 74
$total.Range("D7").Value = 81.15

